# Updates cryptos list prices / 1h-volume changes, and reorders a few rows
# whose coins swapped ranking position, per the commit
# "Updated cryptos list on Sun Dec 10 14:36:55 UTC 2023 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so numeric-looking values
# (e.g. "239.93") are not coerced into Number cells - they must stay text,
# matching the original inlineStr cells.
$ws.Range("D2:D51").NumberFormat = "@"

$updates = @(
    @{ Row=2;  Price="43.777.71" },
    @{ Row=3;  Price="2.346.74";  Change="  -0.80%  " },
    @{ Row=4;  Change="  +0.04%  " },
    @{ Row=5;  Price="239.93";    Change="  -0.65%  " },
    @{ Row=6;  Price="0.665";     Change="  -4.41%  " },
    @{ Row=7;  Price="73.33";     Change="  -4.43%  " },
    @{ Row=8;  Change="  -0.03%  " },
    @{ Row=9;  Change="  -4.56%  " },
    @{ Row=10; Change="  -0.91%  " },
    @{ Row=11; Price="59.70";     Change="  +3.74%  " },
    @{ Row=12; Price="32.69";     Change="  -2.49%  " },
    @{ Row=13; Change="  -0.24%  " },
    @{ Row=14; Price="7.23";      Change="  -3.42%  " },
    @{ Row=15; Price="2.696.22";  Change="  -0.86%  " },
    @{ Row=16; Price="16.08";     Change="  -3.83%  " },
    @{ Row=17; Price="0.902";     Change="  -2.72%  " },
    @{ Row=18; Price="2.345.02";  Change="  -0.35%  " },
    @{ Row=19; Price="43.762.18"; Change="  -0.53%  " },
    @{ Row=20; Price="0.0000103" },
    @{ Row=21; Price="78.54";     Change="  +0.98%  " },
    @{ Row=22; Price="6.64";      Change="  -1.18%  " },
    @{ Row=23; Price="253.04";    Change="  -2.52%  " },

    # Rows 24/25 swap ranking position: WEMIXToken <-> Dai
    @{ Row=24; Coin="Dai"; Link="https://coinranking.com/coin/MoTuySvg7+dai-dai"; Price="1.00"; Change="  -0.13%  " },
    @{ Row=25; Coin="WEMIXToken"; Link="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; Price="3.79"; Change="  +1.60%  " },

    @{ Row=26; Price="1.83"; Change="  +1.79%  " },
    @{ Row=27; Change="  -1.77%  " },

    # Rows 28/29 swap ranking position: Toncoin <-> Cosmos
    @{ Row=28; Coin="Cosmos"; Link="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; Price="10.40"; Change="  -4.52%  " },
    @{ Row=29; Coin="Toncoin"; Link="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; Price="2.31"; Change="  +0.56%  " },

    @{ Row=30; Price="175.19"; Change="  -0.06%  " },
    @{ Row=31; Price="22.18";  Change="  -4.44%  " },
    @{ Row=32; Change="  -0.69%  " },
    @{ Row=33; Change="  -2.56%  " },
    @{ Row=34; Change="  -2.82%  " },
    @{ Row=35; Price="5.06"; Change="  -6.34%  " },
    @{ Row=36; Price="5.32"; Change="  -2.02%  " },
    @{ Row=37; Price="3.80"; Change="  +0.87%  " },
    @{ Row=38; Price="6.37"; Change="  -1.08%  " },
    @{ Row=39; Change="  -2.38%  " },
    @{ Row=40; Price="5.64";   Change="  +16.61%  " },
    @{ Row=41; Price="0.0270"; Change="  -4.45%  " },
    @{ Row=42; Price="64.70";  Change="  +15.64%  " },
    @{ Row=43; Price="9.16";   Change="  -0.59%  " },

    # Rows 44/45 swap ranking position: Cronos <-> InjectiveProtocol
    @{ Row=44; Coin="InjectiveProtocol"; Link="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; Price="18.75"; Change="  -2.99%  " },
    @{ Row=45; Coin="Cronos"; Link="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; Price="0.105"; Change="  -2.94%  " },

    @{ Row=46; Price="0.196"; Change="  -10.73%  " },
    @{ Row=47; Change="  +0.06%  " },
    @{ Row=48; Change="  -3.25%  " },
    @{ Row=49; Change="  -4.18%  " },

    # Rows 50/51 swap ranking position: NEARProtocol <-> Aave
    @{ Row=50; Coin="Aave"; Link="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; Price="98.11"; Change="  -4.54%  " },
    @{ Row=51; Coin="NEARProtocol"; Link="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; Price="2.40"; Change="  -5.75%  " }
)

foreach ($u in $updates) {
    if ($u.Coin) { $ws.Cells.Item($u.Row, 2).Value = $u.Coin }
    if ($u.Link) { $ws.Cells.Item($u.Row, 3).Value = $u.Link }
    if ($u.Price) { $ws.Cells.Item($u.Row, 4).Value = $u.Price }
    if ($u.Change) { $ws.Cells.Item($u.Row, 5).Value = $u.Change }
}
